# Applies the "Excel files and Testfile Bug done Commit" edit:
#  - Merchant(P)CustomFee_QPayNow: Fee Type (I2) and DisplayLabel (M2)
#    corrected from Surcharge/SurchargeFee to CustomFee; view scrolled
#    back to the left with H1 selected.
#  - QPayWithSendlinkCustomFee: view reset to default (no special
#    scroll position / selection).
#  - RunManager: fix the casing of the CustomFeeTxWithQPayNow test name
#    in A18, and select C18.

$wb = $excel.ActiveWorkbook

# --- Sheet: Merchant(P)CustomFee_QPayNow ---
$wsCustomFee = $wb.Worksheets.Item("Merchant(P)CustomFee_QPayNow")
$wsCustomFee.Range("I2").Value = "CustomFee"
$wsCustomFee.Range("M2").Value = "CustomFee"

$wsCustomFee.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$wsCustomFee.Range("H1").Select()

# --- Sheet: QPayWithSendlinkCustomFee ---
$wsSendLink = $wb.Worksheets.Item("QPayWithSendlinkCustomFee")
$wsSendLink.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$wsSendLink.Range("A1").Select()

# --- Sheet: RunManager ---
$wsRunManager = $wb.Worksheets.Item("RunManager")
$wsRunManager.Range("A18").Value = "customFeeTxWithQPayNow_MerchantPrimaryLogin"

$wsRunManager.Activate()
$wsRunManager.Range("C18").Select()
